## "test cases for manage staff"
## - visitorManage: replace the test login row values (recep/dnjnf -> admin,
##   123 -> admin/"123") used to sanity-check the staff-management login flow
## - leaves the selection on visitorManage at B4, and restores SystemUser as
##   the active/selected sheet (C13) once done, matching the saved view state

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("visitorManage")
$wsUsers = $wb.Worksheets.Item("SystemUser")

# Update the two test-login rows to use admin/123 credentials
$wsLogin.Range("A3").Value = "admin"
$wsLogin.Range("B3").Value = "123"
$wsLogin.Range("A4").Value = "admin"
$wsLogin.Range("B4").Value = "123"

# Leave the cursor on B4 of visitorManage (matches the saved selection state)
$wsLogin.Activate() | Out-Null
$wsLogin.Range("B4").Select() | Out-Null

# Re-activate SystemUser (it was the selected tab before the edit) and park
# the cursor at C13, matching the saved selection state
$wsUsers.Activate() | Out-Null
$wsUsers.Range("C13").Select() | Out-Null
